$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the formatting of the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# I0 / IF values for rows 2-70 (row, I-value, J-value)
$ijData = @(
    @(2, 9, 9),
    @(3, 9, 9),
    @(4, 10, 10),
    @(5, 9, 9),
    @(6, 9, 9),
    @(7, 8, 9),
    @(8, 9, 9),
    @(9, 9, 9),
    @(10, 9, 9),
    @(11, 8, 8),
    @(12, 9, 9),
    @(13, 8, 8),
    @(14, 8, 8),
    @(15, 9, 9),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 9, 9),
    @(19, 9, 9),
    @(20, 8, 8),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 9, 9),
    @(25, 8, 8),
    @(26, 9, 9),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 7, 8),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 9, 9),
    @(33, 8, 8),
    @(34, 7, 7),
    @(35, 7, 8),
    @(36, 7, 8),
    @(37, 9, 9),
    @(38, 9, 9),
    @(39, 8, 8),
    @(40, 9, 9),
    @(41, 8, 8),
    @(42, 8, 8),
    @(43, 7, 8),
    @(44, 8, 8),
    @(45, 8, 8),
    @(46, 9, 9),
    @(47, 7, 7),
    @(48, 8, 8),
    @(49, 8, 8),
    @(50, 8, 8),
    @(51, 9, 9),
    @(52, 7, 7),
    @(53, 8, 8),
    @(54, 7, 7),
    @(55, 9, 9),
    @(56, 8, 8),
    @(57, 8, 8),
    @(58, 9, 9),
    @(59, 9, 9),
    @(60, 9, 9),
    @(61, 8, 8),
    @(62, 9, 9),
    @(63, 9, 9),
    @(64, 7, 7),
    @(65, 9, 9),
    @(66, 8, 8),
    @(67, 8, 8),
    @(68, 5, 5),
    @(69, 4, 4),
    @(70, 3, 3)
)

foreach ($entry in $ijData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}
